# "Refactor and user trigger"
#
# The document had several paragraphs whose run-text was split across
# multiple <w:r> elements around <w:proofErr w:type="gramStart"/> /
# <w:proofErr w:type="gramEnd"/> markers (left behind by Word's grammar
# checker). The edit collapses each of those paragraphs back into a single
# run (same visible text) and drops the now-redundant proofErr markers.
# It also relocates the stray "_GoBack" bookmark from the end of the
# "Dalsi formular bude na vyhledavani" paragraph to the very start of the
# "Aplikace bude umoznovat ulozit vybrany binarni obsah..." paragraph
# (the trigger requirement that was moved/renumbered).

$d = $word.ActiveDocument

function Rebuild-Paragraph($index, $text) {
    # Re-writes paragraph $index's content as a single run containing $text,
    # which collapses any run-splits and drops any <w:proofErr/> markers
    # Word's grammar checker had left in the middle of the paragraph.
    #
    # A plain "replace the whole paragraph range" leaves behind a trailing
    # <w:proofErr w:type="gramEnd"/> when that marker is the very last
    # child of the paragraph (nothing textual follows it). Appending a
    # one-character sentinel after the existing content first guarantees
    # the proofErr marker is never the paragraph's last child, so the
    # follow-up full-range rewrite cleanly absorbs (and removes) it.
    $p = $d.Paragraphs.Item($index)
    $rng = $p.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $endPos = $rng.End
    $sentinel = $d.Range($endPos, $endPos)
    $sentinel.InsertAfter("X")

    $p2 = $d.Paragraphs.Item($index)
    $rng2 = $p2.Range
    $rng2.MoveEnd(1, -1) | Out-Null
    $rng2.Text = $text
}

Rebuild-Paragraph 8 "Eviduje a spravuje údaje o jednotlivých uživatelích typu student (jméno, příjmení, kontakty,rok studia, obor, apod.) a o uživatelích typu vyučující (jméno, příjmení, kontakty, vyučované předměty, apod.)"

Rebuild-Paragraph 34 "Administrátor má práva nad vším, nemusí se přepínat. Je to tak OK ?"

Rebuild-Paragraph 40 "Co to jsou číselníky ?"

Rebuild-Paragraph 53 "Stačí to ?"

Rebuild-Paragraph 57 "Funkce na vyhledávání uživatele, vrátí ID ?"

Rebuild-Paragraph 60 "Jsou stejné ? Bude to stačit ?"

Rebuild-Paragraph 63 "Aplikace bude umožňovat uložit vybraný binární obsah do databáze a následně jej i z databáze získat (a pokud se bude jednat o obrázek, tak i v rámci aplikace zobrazit). Pro tento úkol vytvořte ve svém schématu speciální tabulku. Tabulku navrhněte tak, aby kromě samotného binární obsahu umožnila uložit doplňkové informace, jako např.: název souboru, typ souboru, přípona souboru, datum nahrání, datum modifikace, kdo provedl jakou operaci."

Rebuild-Paragraph 68 "Aplikace bude využívat minimálně 3 plnohodnotné formuláře (např. ošetření vstupních polí, apod.) pro vytváření nebo modifikaci dat v tabulkách, ostatní potřebné formuláře jsou samozřejmostí."

# Move the "_GoBack" bookmark from the end of paragraph 71 ("Další
# formulář bude na vyhledávání") to the start of paragraph 63 ("Aplikace
# bude umožňovat uložit vybraný binární obsah...").
if ($d.Bookmarks.Exists("_GoBack")) {
    $old = $d.Bookmarks.Item("_GoBack")
    $old.Delete()
}
$target = $d.Paragraphs.Item(63).Range
$targetStart = $d.Range($target.Start, $target.Start)
$d.Bookmarks.Add("_GoBack", $targetStart) | Out-Null
